# Rename the stock ticker "APPL" -> "AAPL" everywhere it appears in the
# workbook (it was a typo: "APPL should always have been AAPL").
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "APPL") {
            $cell.Value2 = "AAPL"
        }
    }
}
